$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.636.44'
$ws.Range("E2").Value = '  -5.83%  '
$ws.Range("D3").Value = '3.333.59'
$ws.Range("E3").Value = '  -5.03%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.64%  '
$ws.Range("E7").Value = '  -9.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.382'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -10.97%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.923'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -12.69%  '
$ws.Range("D11").Value = '3.338.76'
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.193'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -13.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.41%  '
$ws.Range("D15").Value = '92.472.48'
$ws.Range("E15").Value = '  -5.87%  '
$ws.Range("D16").Value = '3.957.35'
$ws.Range("E16").Value = '  -5.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000242'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -12.81%  '
$ws.Range("D19").Value = '3.330.62'
$ws.Range("E19").Value = '  -5.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -12.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.446'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -15.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -10.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000183'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.48%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.523.81'
$ws.Range("E28").Value = '  -4.69%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.17%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.60%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.96%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.130'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -11.39%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.992'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.169'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -11.78%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.53%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.518'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -12.37%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.53%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '513.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.30%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.146'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.23%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.863'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.40%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.74%  '
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.93%  '
$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.55'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.32%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.79%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0387'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -10.96%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.03%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.27%  '
